# RoboRTS-Firmware Parameters.xlsx — "Add files via upload / PID debugging 3/17"
#
# Adds a new "Gimbal_PID" worksheet (after the existing "sys_config.h" sheet)
# containing three dated snapshots (3/5, 3/10, 3/17 2019) of the gimbal /
# chassis-speed PID tuning parameters, and leaves the active selection on the
# new sheet (matching the author's saved view).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- park the selection on the old sheet where the diff shows it ending up ---
$ws1.Range("C4").Select()

# --- insert the new sheet right after sys_config.h ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Gimbal_PID"

# ---------------------------------------------------------------------------
# Header row (1): "Parameter" + the 12 PID/speed-loop column names
# ---------------------------------------------------------------------------
$ws2.Range("A1").Value = "Parameter"
$ws2.Range("A1").HorizontalAlignment = -4108   # xlCenter
$ws2.Range("A1").VerticalAlignment = -4160     # xlTop
$ws2.Range("A1").WrapText = $true

$headers = @("kp_pit","ki_pit","kd_pit","kp_yaw","ki_yaw","kd_yaw", `
             "speed_p_kp","speed_p_ki","speed_p_kd","speed_y_kp","speed_y_ki","speed_y_kd")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws2.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# Row 2: "Description" label only (no description text filled in yet)
$ws2.Range("A2").Value = "Description"
$ws2.Range("A2").HorizontalAlignment = -4108   # xlCenter
$ws2.Range("A2").VerticalAlignment = -4160     # xlTop
$ws2.Range("A2").WrapText = $true

# ---------------------------------------------------------------------------
# Data rows: date in column A, 12 numeric values in B:M
# ---------------------------------------------------------------------------
$rows = @(
    @{ Date = 43529; Values = @(156, 0,   205.2, 156,  1,  450, 3, 0, 3, 3,  0,  1.2) },  # 3/5/2019
    @{ Date = 43534; Values = @(170, 3.5, 100,   1200, 15, 30,  3, 0, 3, 10, 0,  50)  },  # 3/10/2019
    @{ Date = 43541; Values = @(170, 3.5, 100,   650,  0,  150, 3, 0, 3, 20, -8, 35)  }   # 3/17/2019
)

$r = 3
foreach ($entry in $rows) {
    $cell = $ws2.Cells.Item($r, 1)
    $cell.Value = $entry.Date
    $cell.NumberFormat = "m/d/yyyy"
    $cell.HorizontalAlignment = -4131   # xlLeft
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.WrapText = ($r -eq 3)         # only the first (template) row wraps

    $vals = $entry.Values
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws2.Cells.Item($r, $i + 2).Value = $vals[$i]
    }
    $r++
}

# ---------------------------------------------------------------------------
# Column widths
# ---------------------------------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 13.25
$ws2.Columns.Item(8).ColumnWidth = 15.83
$ws2.Columns.Item(9).ColumnWidth = 12.83
$ws2.Columns.Item(10).ColumnWidth = 15
$ws2.Columns.Item(11).ColumnWidth = 12.83
$ws2.Columns.Item(12).ColumnWidth = 12
$ws2.Columns.Item(13).ColumnWidth = 13.58

# ---------------------------------------------------------------------------
# "Mark changed value" conditional formatting, same pattern already used on
# sys_config.h: highlight a row when it differs from the row above it.
# ---------------------------------------------------------------------------
$fc = $ws2.Range("B4:XFD4").FormatConditions.Add(1, 4, "=B3")
$fc.Interior.Color = 5287936
$fc = $ws2.Range("A4").FormatConditions.Add(1, 4, "=A3")
$fc.Interior.Color = 5287936
$fc = $ws2.Range("B5:XFD5").FormatConditions.Add(1, 4, "=B4")
$fc.Interior.Color = 5287936
$fc = $ws2.Range("A5").FormatConditions.Add(1, 4, "=A4")
$fc.Interior.Color = 5287936

# ---------------------------------------------------------------------------
# Make the new sheet the active / visible one, with the same selection the
# author left it on.
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("B5").Select()
